$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural change: the "Quantity" column (F) moves to become column B.
# Cutting column F and inserting it before column B shifts Value/Device/Package/
# Description from B:E to C:F, and Excel automatically repoints the LINE_PRICE
# formulas (=F*G) to the new Quantity column (=B*G).
$ws.Columns("F:F").Cut() | Out-Null
$ws.Columns("B:B").Insert() | Out-Null

# --- BOM content tweaks (low-stock substitutions etc.) ---

# Row 5: 100nF ceramic caps - quantity trimmed from 11 to 10
$ws.Range("B5").Value = 10

# Row 19: Schottky diode rows now carry an explicit "Value"
$ws.Range("C19").Value = "SCHOTTKY"

# Row 22: LM1117 3.3V regulator replaced with TLV1117 (low stock substitute)
$ws.Range("C22").Value = "TLV1117LV33DCYT"
$ws.Range("F22").Value = "Voltage Regulator TLV1117. 1 Amp. Pos."
$ws.Range("G22").Value = 1.1399999999999999
$ws.Range("J22").Value = "TLV1117LV33DCYT"
$ws.Range("K22").Value = "595-TLV1117LV33DCYT"
$ws.Range("L22").Value = "https://www.mouser.com/ProductDetail/Texas-Instruments/TLV1117LV33DCYT?qs=tEz3BkPb1ry8fvdWR4Nyog%3D%3D"

# Row 28: ESP32-WROOM-32U swapped for ESP32-WROOM-32UE-N4 (low stock substitute)
$ws.Range("C28").Value = "ESP32-WROOM-32UE-N4"
$ws.Range("D28").Value = "ESP32-WROOM-32UE-N4"
$ws.Range("F28").Value = "ESP32-D0WD, 4Mbits SPI flash, U.FL antenna"
$ws.Range("G28").Value = 3
$ws.Range("J28").Value = "ESP32-WROOM-32UE-N4"
$ws.Range("K28").Value = "356-ESP32WRM32UE32UH"

# Row 40: CH340G sourcing note expanded
$ws.Range("L40").Value = "Not commonly available from western distributors. May need to buy on eBay, Amazon, AliExpress. Sparkfun sell the CH340E (different package), but only in quantity 10."

# --- Cosmetic view changes ---
$ws.Range("A1").Select() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("H46").Select() | Out-Null

# --- Column width tweaks ---
$ws.Columns("A:A").ColumnWidth = 38.69140625
$ws.Columns("B:B").ColumnWidth = 10.07421875
